# Historias de Usuario (Clientes): rename the lone worksheet from the
# generic default "Hoja1" to "Cliente", and leave the selection on O12
# (where the author last clicked before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Cliente"

[void]$ws.Range("O12").Select()
